$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '27.824.11'
Set-TextValue 'E2' '  +1.68%  '
Set-TextValue 'D3' '1.887.12'
Set-TextValue 'E3' '  +1.57%  '
Set-TextValue 'D4' '1.007'
Set-TextValue 'E4' '  +0.64%  '
Set-TextValue 'D5' '334.37'
Set-TextValue 'E5' '  +1.21%  '
Set-TextValue 'D6' '1.007'
Set-TextValue 'E6' '  +0.50%  '
Set-TextValue 'D7' '0.4719'
Set-TextValue 'E7' '  +1.35%  '
Set-TextValue 'D8' '0.3928'
Set-TextValue 'E8' '  -0.62%  '
Set-TextValue 'D9' '47.75'
Set-TextValue 'E9' '  +1.73%  '
Set-TextValue 'D10' '0.08062'
Set-TextValue 'E10' '  +1.00%  '
Set-TextValue 'D11' '1.027'
Set-TextValue 'E11' '  +1.07%  '
Set-TextValue 'D12' '22.12'
Set-TextValue 'E12' '  +2.66%  '
Set-TextValue 'D13' '1.876.14'
Set-TextValue 'E13' '  +1.66%  '
Set-TextValue 'D14' '5.991'
Set-TextValue 'E14' '  +0.75%  '
Set-TextValue 'E15' '  -0.10%  '
Set-TextValue 'D16' '1.010'
Set-TextValue 'E16' '  +0.93%  '
Set-TextValue 'D17' '0.06729'
Set-TextValue 'E17' '  +2.58%  '
Set-TextValue 'B18' 'Litecoin'
Set-TextValue 'C18' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D18' '87.30'
Set-TextValue 'E18' '  +1.12%  '
Set-TextValue 'B19' 'ShibaInu'
Set-TextValue 'C19' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D19' '0.00001048'
Set-TextValue 'E19' '  +1.26%  '
Set-TextValue 'D20' '17.32'
Set-TextValue 'E20' '  +0.13%  '
Set-TextValue 'E21' '  +0.48%  '
Set-TextValue 'D22' '27.856.70'
Set-TextValue 'E22' '  +1.82%  '
Set-TextValue 'D23' '5.524'
Set-TextValue 'E23' '  +0.53%  '
Set-TextValue 'D24' '10.98'
Set-TextValue 'E24' '  +0.57%  '
Set-TextValue 'D25' '2.331'
Set-TextValue 'E25' '  +1.38%  '
Set-TextValue 'D26' '2.104.46'
Set-TextValue 'E26' '  +1.84%  '
Set-TextValue 'D27' '158.99'
Set-TextValue 'E27' '  +3.42%  '
Set-TextValue 'D28' '20.13'
Set-TextValue 'E28' '  -1.82%  '
Set-TextValue 'D29' '2.105'
Set-TextValue 'E29' '  +1.68%  '
Set-TextValue 'E30' '  +1.65%  '
Set-TextValue 'D31' '122.06'
Set-TextValue 'E31' '  -0.09%  '
Set-TextValue 'D32' '0.9785'
Set-TextValue 'E32' '  +2.74%  '
Set-TextValue 'D33' '0.09504'
Set-TextValue 'E33' '  +0.19%  '
Set-TextValue 'D34' '1.451'
Set-TextValue 'E34' '  +0.23%  '
Set-TextValue 'D35' '3.622'
Set-TextValue 'E35' '  +1.04%  '
Set-TextValue 'D36' '5.357'
Set-TextValue 'E36' '  +1.53%  '
Set-TextValue 'D37' '0.06166'
Set-TextValue 'E37' '  +1.74%  '
Set-TextValue 'D38' '0.02272'
Set-TextValue 'E38' '  +1.69%  '
Set-TextValue 'D39' '1.218'
Set-TextValue 'E39' '  +0.35%  '
Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '8.069'
Set-TextValue 'E40' '  +0.19%  '
Set-TextValue 'B41' 'TheSandbox'
Set-TextValue 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D41' '0.6009'
Set-TextValue 'E41' '  +1.00%  '
Set-TextValue 'E42' '  +0.06%  '
Set-TextValue 'E43' '  +0.50%  '
Set-TextValue 'D44' '1.260'
Set-TextValue 'E44' '  -0.71%  '
Set-TextValue 'D45' '0.5711'
Set-TextValue 'E45' '  +0.99%  '
Set-TextValue 'D46' '12.20'
Set-TextValue 'E46' '  +0.83%  '
Set-TextValue 'D47' '3.401'
Set-TextValue 'E47' '  -0.99%  '
Set-TextValue 'D48' '1.942'
Set-TextValue 'E48' '  +0.63%  '
Set-TextValue 'D49' '0.06915'
Set-TextValue 'E49' '  +2.16%  '
Set-TextValue 'D50' '113.29'
Set-TextValue 'E50' '  +3.47%  '
Set-TextValue 'D51' '0.00000000301'
Set-TextValue 'E51' '  +6.57%  '
